$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 58; this pushes old rows 58-128 down to 59-129
$ws.Rows.Item(58).Insert()

# Populate the new row 58 with the new record data
$ws.Cells.Item(58, 1).Value = 10
$ws.Cells.Item(58, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(58, 3).Value = "La Araucanía"
$ws.Cells.Item(58, 4).Value = 44413
$ws.Cells.Item(58, 5).Value = 9
$ws.Cells.Item(58, 6).Value = 100112043
$ws.Cells.Item(58, 7).Value = "Pepino dulce"
$ws.Cells.Item(58, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(58, 9).Value = "Primera"
$ws.Cells.Item(58, 10).Value = 250
$ws.Cells.Item(58, 11).Value = 19000
$ws.Cells.Item(58, 12).Value = 20000
$ws.Cells.Item(58, 13).Value = 19600
$ws.Cells.Item(58, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(58, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(58, 16).Value = 1089
$ws.Cells.Item(58, 17).Value = 18
$ws.Cells.Item(58, 18).Value = "Hortaliza"
